# LOQ4204.xlsx update
# - "Objetivos:" body (B10/C10) replaced with the professor's name
# - A new "Programa resumido:" row (row 13) is inserted with body "Semestral"
# - The old "Programa resumido:"/"Programa:" bodies are dropped, and the
#   remaining labels in column A (rows 14-21) shift up by one
# - The long Portuguese "Programa:" paragraph (old row 16 body) is removed
# - Row 17 ("Avaliação:") loses its B/C body cells and its custom row height
# - A new body (the professor's name again) is added at B18/C18 under "Método:"
# - The final row (old row 22, "Bibliografia:" + long bibliography text) is deleted

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats / xlPasteValues constants (used below to move values between
# cells while picking up the destination column's existing number format /
# style, instead of Excel's "Value = ..." auto-guessing e.g. dates)
$xlPasteFormats = -4122
$xlPasteValues = -4163

# Objetivos: body -> professor name
$ws.Range("B10").Value = "11079086 - Herlandí de Souza Andrade"
$ws.Range("C10").Value = "11079086 - Herlandí de Souza Andrade"

# New row 13: "Programa resumido:" / "Semestral"
$ws.Rows(13).RowHeight = 60
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Row 14: label shifts to "Short syllabus:", body becomes the short syllabus (EN)
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = "A. Microeconomy. B. Macroeconomy. C. Economic Development. D. International Economy. E. Brazilian Economy"
$ws.Range("C14").Value = "A. Microeconomy. B. Macroeconomy. C. Economic Development. D. International Economy. E. Brazilian Economy"

# Row 15: label shifts to "Programa:"; body becomes (erroneously) the activation
# date, copied in as a value so it stays text instead of becoming a date serial
$ws.Rows(15).RowHeight = 120
$ws.Range("A15").Value = "Programa:"
$ws.Range("B8").Copy()
$ws.Range("B15").PasteSpecial($xlPasteValues)
$ws.Range("C8").Copy()
$ws.Range("C15").PasteSpecial($xlPasteValues)

# Row 16: label shifts to "Syllabus:"; body keeps the long English syllabus text
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = "A. MICROECONOMY: 1. Introduction to the concepts of Economics and fundamentals of microeconomic analysis. 2. Consumer and demand theory. 3. Firm and offer theory. 4. Costs and price formation. 5. Market Structures 6. Strategic behavior and competition. 7. Technology as a factor of production. 8. Sustainability: resources, costs and environmental indicators.B. MACROECONOMY: 1. Fundamentals of macroeconomic analysis. 2. National accounting. 3. Classical and Keynesian balances. 4. Monetary system. 5. Fiscal policy. 6. World economy and international trade. 7. Fundamentals of regression as a tool to quantify economic relationships. 8. Public sector.C. ECONOMIC DEVELOPMENT: 1. Growth factors. 2. Sources of Development. 3. Financing Economic Development. 4. A model of economic growth. 5. The internationalization and globalization process.D. INTERNATIONAL ECONOMY: 1. Fundamentals of International Trade. 2. Determination of Exchange Rates. 3. External policies. 4. Factors determining the behavior of imports and exports.E. BRAZILIAN ECONOMY: 1. The historical experience of Brazilian industrialization. 2. The internationalization of the Brazilian economy. 3. Cycle theory and Brazilian reality. 4. Brazil's economic cycles throughout its recent history."
$ws.Range("C16").Value = "A. MICROECONOMY: 1. Introduction to the concepts of Economics and fundamentals of microeconomic analysis. 2. Consumer and demand theory. 3. Firm and offer theory. 4. Costs and price formation. 5. Market Structures 6. Strategic behavior and competition. 7. Technology as a factor of production. 8. Sustainability: resources, costs and environmental indicators.B. MACROECONOMY: 1. Fundamentals of macroeconomic analysis. 2. National accounting. 3. Classical and Keynesian balances. 4. Monetary system. 5. Fiscal policy. 6. World economy and international trade. 7. Fundamentals of regression as a tool to quantify economic relationships. 8. Public sector.C. ECONOMIC DEVELOPMENT: 1. Growth factors. 2. Sources of Development. 3. Financing Economic Development. 4. A model of economic growth. 5. The internationalization and globalization process.D. INTERNATIONAL ECONOMY: 1. Fundamentals of International Trade. 2. Determination of Exchange Rates. 3. External policies. 4. Factors determining the behavior of imports and exports.E. BRAZILIAN ECONOMY: 1. The historical experience of Brazilian industrialization. 2. The internationalization of the Brazilian economy. 3. Cycle theory and Brazilian reality. 4. Brazil's economic cycles throughout its recent history."

# Row 17: label shifts to "Avaliação:"; body cells and the custom row height
# are both removed (back to the sheet's default row height)
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17").Clear()
$ws.Range("C17").Clear()
$ws.Rows(17).AutoFit()

# Row 18: label shifts to "Método:"; body becomes the professor's name again.
# Paste formats first (so B18/C18 pick up the normal column style) then paste
# the value on top of it.
$ws.Rows(18).RowHeight = 60
$ws.Range("A18").Value = "Método:"
$ws.Range("B14").Copy()
$ws.Range("B18").PasteSpecial($xlPasteFormats)
$ws.Range("B10").Copy()
$ws.Range("B18").PasteSpecial($xlPasteValues)
$ws.Range("C14").Copy()
$ws.Range("C18").PasteSpecial($xlPasteFormats)
$ws.Range("C10").Copy()
$ws.Range("C18").PasteSpecial($xlPasteValues)

# Row 19: label shifts to "Critério:" (body - the teaching method text - stays the same)
$ws.Range("A19").Value = "Critério:"

# Row 20: label shifts to "Norma de recuperação:" (body - the evaluation criteria - stays the same)
$ws.Range("A20").Value = "Norma de recuperação:"

# Row 21: label shifts to "Bibliografia:"; body stays the recovery-exam formula, height grows to 120
$ws.Rows(21).RowHeight = 120
$ws.Range("A21").Value = "Bibliografia:"

# Row 22 (old "Bibliografia:" label + the long bibliography paragraph) is dropped entirely
$ws.Rows(22).Delete()
